# Weekly roll of the "Hortaliza, Femacal de La Calera - Ciboulette" price
# series: insert one new week of data at the top of the data block (row 77,
# the first data row for this subset) and push the existing rows down by
# one, dropping the oldest row off the bottom into a newly appended row 311.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 77:310 down to 78:311, preserving formatting/styles, and
# leave a blank row 77 ready for the new week's data.
$ws.Rows.Item(77).Insert()

# Populate the new week's row with the latest observation.
$ws.Range("A77").Value = 3
$ws.Range("B77").Value = "Femacal de La Calera"
$ws.Range("C77").Value = "Coquimbo"
$ws.Range("D77").Value = 44690
$ws.Range("E77").Value = 5
$ws.Range("F77").Value = 100112039
$ws.Range("G77").Value = "Ciboulette"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 160
$ws.Range("K77").Value = 1500
$ws.Range("L77").Value = 1500
$ws.Range("M77").Value = 1500
$ws.Range("N77").Value = '$/docena de atados'
$ws.Range("O77").Value = "Provincia de Quillota"
$ws.Range("P77").Value = 500
$ws.Range("Q77").Value = 3
$ws.Range("R77").Value = "Hortaliza"
